$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2916
$ws.Range("E2").Value = 130
$ws.Range("F2").Value = 130
$ws.Range("G2").Value = 120
$ws.Range("H2").Value = 93
$ws.Range("I2").Value = 93
$ws.Range("K2").Value = 1886
$ws.Range("L2").Value = 689
$ws.Range("M2").Value = 1197
$ws.Range("N2").Value = 1197
$ws.Range("P2").Value = 174
$ws.Range("Q2").Value = 178
$ws.Range("R2").Value = -15
$ws.Range("S2").Value = -133
$ws.Range("T2").Value = 14
$ws.Range("U2").Value = 165
$ws.Range("V2").Value = 398
$ws.Range("W2").Value = 4.46
$ws.Range("X2").Value = 3.19
$ws.Range("AA2").Value = 57.58
$ws.Range("AB2").Value = 585.74
$ws.Range("AC2").Value = 267
$ws.Range("AD2").Value = 9.02
$ws.Range("AE2").Value = 3432
$ws.Range("AF2").Value = 0.7
$ws.Range("AG2").Value = 70
$ws.Range("AH2").Value = 2.91
$ws.Range("AI2").Value = 26.25
$ws.Range("AJ2").Value = 34869420
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()

# Row 3
$ws.Range("D3").Value = 2620
$ws.Range("E3").Value = 181
$ws.Range("F3").Value = 181
$ws.Range("G3").Value = 194
$ws.Range("H3").Value = 145
$ws.Range("I3").Value = 145
$ws.Range("K3").Value = 1875
$ws.Range("L3").Value = 570
$ws.Range("M3").Value = 1305
$ws.Range("N3").Value = 1305
$ws.Range("P3").Value = 174
$ws.Range("Q3").Value = 175
$ws.Range("R3").Value = -90
$ws.Range("S3").Value = -167
$ws.Range("T3").Value = 18
$ws.Range("U3").Value = 157
$ws.Range("V3").Value = 262
$ws.Range("W3").Value = 6.9
$ws.Range("X3").Value = 5.53
$ws.Range("Y3").Value = 11.58
$ws.Range("Z3").Value = 7.7
$ws.Range("AA3").Value = 43.64
$ws.Range("AB3").Value = 652.99
$ws.Range("AC3").Value = 415
$ws.Range("AD3").Value = 7.89
$ws.Range("AE3").Value = 3743
$ws.Range("AF3").Value = 0.88
$ws.Range("AG3").Value = 100
$ws.Range("AH3").Value = 3.05
$ws.Range("AI3").Value = 24.07
$ws.Range("AJ3").Value = 34869420
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 2312
$ws.Range("E4").Value = 166
$ws.Range("F4").Value = 181
$ws.Range("G4").Value = 169
$ws.Range("H4").Value = 140
$ws.Range("I4").Value = 140
$ws.Range("K4").Value = 1917
$ws.Range("L4").Value = 518
$ws.Range("M4").Value = 1399
$ws.Range("N4").Value = 1399
$ws.Range("P4").Value = 174
$ws.Range("Q4").Value = 296
$ws.Range("R4").Value = -22
$ws.Range("S4").Value = -73
$ws.Range("T4").Value = 33
$ws.Range("U4").Value = 263
$ws.Range("V4").Value = 227
$ws.Range("W4").Value = 7.2
$ws.Range("X4").Value = 6.05
$ws.Range("Y4").Value = 10.35
$ws.Range("Z4").Value = 7.38
$ws.Range("AA4").Value = 37.06
$ws.Range("AB4").Value = 712.18
$ws.Range("AC4").Value = 401
$ws.Range("AD4").Value = 8.15
$ws.Range("AE4").Value = 4012
$ws.Range("AF4").Value = 0.82
$ws.Range("AG4").Value = 115
$ws.Range("AH4").Value = 3.52
$ws.Range("AI4").Value = 28.67
$ws.Range("AJ4").Value = 34869420
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 4901
$ws.Range("E5").Value = 136
$ws.Range("F5").Value = 136
$ws.Range("G5").Value = 104
$ws.Range("H5").Value = 151
$ws.Range("I5").Value = 151
$ws.Range("K5").Value = 2685
$ws.Range("L5").Value = 1201
$ws.Range("M5").Value = 1484
$ws.Range("N5").Value = 1484
$ws.Range("P5").Value = 174
$ws.Range("Q5").Value = 41
$ws.Range("R5").Value = -482
$ws.Range("S5").Value = 375
$ws.Range("T5").Value = 41
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 817
$ws.Range("W5").Value = 2.78
$ws.Range("X5").Value = 3.08
$ws.Range("Y5").Value = 10.48
$ws.Range("Z5").Value = 6.56
$ws.Range("AA5").Value = 80.94
$ws.Range("AB5").Value = 777.38
$ws.Range("AC5").Value = 433
$ws.Range("AD5").Value = 7.85
$ws.Range("AE5").Value = 4407
$ws.Range("AF5").Value = 0.77
$ws.Range("AG5").Value = 160
$ws.Range("AH5").Value = 4.71
$ws.Range("AI5").Value = 35.67
$ws.Range("AJ5").Value = 34869420
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 8309
$ws.Range("E6").Value = 186
$ws.Range("F6").Value = 186
$ws.Range("G6").Value = 204
$ws.Range("H6").Value = 163
$ws.Range("I6").Value = 163
$ws.Range("K6").Value = 2801
$ws.Range("L6").Value = 1202
$ws.Range("M6").Value = 1599
$ws.Range("N6").Value = 1599
$ws.Range("P6").Value = 174
$ws.Range("Q6").Value = 124
$ws.Range("R6").Value = -133
$ws.Range("S6").Value = -108
$ws.Range("T6").Value = 95
$ws.Range("U6").Value = 30
$ws.Range("V6").Value = 762
$ws.Range("W6").Value = 2.23
$ws.Range("X6").Value = 1.96
$ws.Range("Y6").Value = 10.57
$ws.Range("Z6").Value = 5.94
$ws.Range("AA6").Value = 75.17
$ws.Range("AB6").Value = 847.93
$ws.Range("AC6").Value = 467
$ws.Range("AD6").Value = 7.05
$ws.Range("AE6").Value = 4749
$ws.Range("AF6").Value = 0.69
$ws.Range("AG6").Value = 160
$ws.Range("AH6").Value = 4.86
$ws.Range("AI6").Value = 33.06
$ws.Range("AJ6").Value = 34869420
$ws.Range("J6").ClearContents()
$ws.Range("O6").ClearContents()

# Row 7
$ws.Range("D7:AJ7").ClearContents()

# Row 8
$ws.Range("D8:AJ8").ClearContents()

# Row 9
$ws.Range("D9:AJ9").ClearContents()
